# Update build/version timestamp strings across the workbook.
# Old build timestamp -> new build timestamp
$oldTimestamp = "February 03 2026 17.29.55 EST"
$newTimestamp = "February 03 2026 18.05.36 EST"

$oldVersionString = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $oldTimestamp)"
$newVersionString = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newTimestamp)"

$wb = $excel.ActiveWorkbook

# --- "About" sheet ---
$aboutSheet = $wb.Worksheets.Item("About")

$aboutSheet.Range("A2").Value = "Version: $newVersionString"

$aboutSheet.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Lu'an Licun Coal Mine, China, M1170, version '$newVersionString'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 9; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # Column S = 19th column (build_version)
    if ($cell.Value2 -eq $oldVersionString) {
        $cell.Value = $newVersionString
    }
}
